$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New glossary rows for the "UST Fair" Future-promo photo/audio set (rows 62-67).
# Columns: A=Key Words_zh-HK, B=Key Words_zh-CN, C=Key Words_en-US,
#          D=Image Name, E=Voice Name.
# "/鱻/" is this sheet's existing placeholder for "no keyword in this language".
$filler = "/鱻/"

# Image Name column (D), entered in file order 1,2,3,5,4,6.
$ws.Cells.Item(62, 4).Value = "Future 1.jpg"
$ws.Cells.Item(63, 4).Value = "Future 2.jpg"
$ws.Cells.Item(64, 4).Value = "Future 3.jpg"
$ws.Cells.Item(66, 4).Value = "Future 5.jpg"
$ws.Cells.Item(65, 4).Value = "Future 4.jpg"
$ws.Cells.Item(67, 4).Value = "Future 6.jpg"

# Voice Name column (E), entered in row order.
$ws.Cells.Item(62, 5).Value = "Section1 Welcome.mp3"
$ws.Cells.Item(63, 5).Value = "Section2 Always.mp3"
$ws.Cells.Item(64, 5).Value = "Section3 Journey.mp3"
$ws.Cells.Item(65, 5).Value = "Section4 Anything.mp3"
$ws.Cells.Item(66, 5).Value = "Section5 Engineering.mp3"
$ws.Cells.Item(67, 5).Value = "Section6 Have.mp3"

# Key Words_en-US column (C), entered out of row order.
$ws.Cells.Item(67, 3).Value = "/good day/,/bye/"
$ws.Cells.Item(64, 3).Value = "/journey/,/destination/"
$ws.Cells.Item(63, 3).Value = "/curious/,/embrace/,/"
$ws.Cells.Item(62, 3).Value = "/welcome/,/nervous/,/interview/"
$ws.Cells.Item(65, 3).Value = "/anything/,/design/"
$ws.Cells.Item(66, 3).Value = "/engineering/"

# Key Words_zh-HK / Key Words_zh-CN columns (A/B) - no localized keyword yet.
foreach ($row in 62..67) {
    $ws.Cells.Item($row, 1).Value = $filler
    $ws.Cells.Item($row, 2).Value = $filler
}

# Move the glossary's visual cursor to where editing left off.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 37
$win.ScrollColumn = 1
$ws.Range("C66").Select()
